$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared string "PENDENTE" by setting it into M2 (replaces FATURADO)
$ws.Range("M2").Value = "PENDENTE"

# Set column M width to (as closely as possible) match the bestFit width
# from the diff (13.7109375 character units)
$ws.Columns.Item(13).ColumnWidth = 12.833333333333334

# Update the selected cell to M3
$ws.Range("M3").Select() | Out-Null
